# Revert dev content to version bf6ab339:
#  - remove the "Address" column (I) entirely, shifting
#    Country / Longitude / Latitude one column to the left (I/J/K)
#  - the old Longitude/Latitude data columns (now J and K) no longer
#    carry real coordinates, so blank them out to "NA"
#  - fix the "Latitude" header typo -> "Lattitude"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column I ("Address"); J->I, K->J, L->K shift left.
$ws.Range("I1").EntireColumn.Delete()

# Former Longitude/Latitude values are no longer tracked -> "NA".
$ws.Range("J2:J26").Value = "NA"
$ws.Range("K2:K26").Value = "NA"

# Header typo fix.
$ws.Range("K1").Value = "Lattitude"
